$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to the batter's name
$ws.Name = "Pat Cummins"

# Headers (row 1) - new column order with matchNo inserted as the first column
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data rows 2-6: matchNo, teamName, batterName, states, runs, balls, fours, sixes, sr, opponentTeamName, venue, date, result
$data = @(
    @("10th","Kolkata Knight Riders","Pat Cummins","c †de Villiers b Jamieson","6","2","0","1","300.00","Royal Challengers Bangalore","Chennai","April 18","RCB won by 38 runs"),
    @("25th","Kolkata Knight Riders","Pat Cummins","","11","13","1","0","84.61","Delhi Capitals","Ahmedabad","April 29","Capitals won by 7 wickets (with 21 balls remaining)"),
    @("15th","Kolkata Knight Riders","Pat Cummins","","66","34","4","6","194.11","Chennai Super Kings","Wankhede","April 21","Super Kings won by 18 runs"),
    @("18th","Kolkata Knight Riders","Pat Cummins","c Parag b Morris","10","6","0","1","166.66","Rajasthan Royals","Wankhede","April 24","Royals won by 6 wickets (with 7 balls remaining)"),
    @("5th","Kolkata Knight Riders","Pat Cummins","b Boult","0","1","0","0","0.00","Mumbai Indians","Chennai","April 13","Mumbai won by 10 runs")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        # The source data stores every field as text. Empty values and
        # numeric-looking values (runs, balls, fours, sixes, sr, ...) need a
        # leading quote so Excel keeps them as text instead of coercing them
        # to a number / blank cell.
        if ($val -eq "" -or $val -match '^-?\d+(\.\d+)?$') {
            $ws.Cells.Item($r + 2, $c + 1).Value = "'" + $val
        } else {
            $ws.Cells.Item($r + 2, $c + 1).Value = $val
        }
    }
}
